$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws2 = $wb.Worksheets.Item("Tabelle2")

# --- Tabelle1: add new data rows for 2018, 2019, 2020 (copy formatting from the 2017 row) ---
$ws1.Range("A12:B12").Copy($ws1.Range("A13:B13"))
$ws1.Range("A13").Value = 2018
$ws1.Range("A12:B12").Copy($ws1.Range("A14:B14"))
$ws1.Range("A14").Value = 2019
$ws1.Range("A12:B12").Copy($ws1.Range("A15:B15"))
$ws1.Range("A15").Value = 2020

# --- Tabelle1: shift the "Datenquelle / Copyright / Vervielfaeltigung" footer block down one row ---
$ws1.Rows("35").Insert(-4121)
$ws1.Range("B36").Value = "Weltgesundheitsorganisation (WHO)"
$ws1.Range("B37").Value = "©       Statistisches Bundesamt (Destatis) 2021"
$ws1.Rows("42").Delete()
$ws1.Rows("41").Delete()

# --- Tabelle2 (chart source data): add rows for 2018, 2019, 2020 ---
$ws2.Range("A8:C8").Copy($ws2.Range("A9:C9"))
$ws2.Range("A9").Value = 2018
$ws2.Range("A8:C8").Copy($ws2.Range("A10:C10"))
$ws2.Range("A10").Value = 2019
$ws2.Range("A8:C8").Copy($ws2.Range("A11:C11"))
$ws2.Range("A11").Value = 2020
